$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.795.55"
$ws.Range("E2").Value = "  +4.33%  "
$ws.Range("D3").Value = "3.440.86"
$ws.Range("E3").Value = "  +3.73%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "579.17"
$ws.Range("E5").Value = "  +4.86%  "
$ws.Range("D6").Value = "184.67"
$ws.Range("E6").Value = "  +6.56%  "
$ws.Range("D7").Value = "0.633"
$ws.Range("E7").Value = "  +2.81%  "
$ws.Range("D8").Value = "3.435.10"
$ws.Range("E8").Value = "  +3.65%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("E10").Value = "  +2.39%  "
$ws.Range("D11").Value = "0.646"
$ws.Range("E11").Value = "  +2.64%  "
$ws.Range("D12").Value = "56.20"
$ws.Range("E12").Value = "  +5.15%  "
$ws.Range("E13").Value = "  +0.33%  "
$ws.Range("D14").Value = "9.43"
$ws.Range("E14").Value = "  +4.58%  "
$ws.Range("D15").Value = "3.992.47"
$ws.Range("E15").Value = "  +3.79%  "
$ws.Range("D16").Value = "18.69"
$ws.Range("E16").Value = "  +3.52%  "
$ws.Range("D17").Value = "3.437.59"
$ws.Range("E17").Value = "  +3.15%  "
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("D19").Value = "66.728.49"
$ws.Range("E19").Value = "  +2.80%  "
$ws.Range("D20").Value = "12.08"
$ws.Range("E20").Value = "  +3.56%  "
$ws.Range("D21").Value = "1.02"
$ws.Range("E21").Value = "  +3.54%  "
$ws.Range("D22").Value = "485.84"
$ws.Range("E22").Value = "  +7.75%  "
$ws.Range("D23").Value = "16.82"
$ws.Range("E23").Value = "  +21.92%  "
$ws.Range("D24").Value = "5.04"
$ws.Range("E24").Value = "  +1.94%  "
$ws.Range("D25").Value = "4.37"
$ws.Range("E25").Value = "  +7.56%  "
$ws.Range("D26").Value = "89.77"
$ws.Range("E26").Value = "  +3.68%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").Value = "11.08"
$ws.Range("E27").Value = "  +4.08%  "
$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D28").Value = "2.95"
$ws.Range("E28").Value = "  +3.38%  "
$ws.Range("D29").Value = "9.16"
$ws.Range("E29").Value = "  +7.10%  "
$ws.Range("D30").Value = "31.36"
$ws.Range("E30").Value = "  +1.74%  "
$ws.Range("D31").Value = "7.09"
$ws.Range("E31").Value = "  +8.54%  "
$ws.Range("D32").Value = "64.54"
$ws.Range("E32").Value = "  +6.88%  "
$ws.Range("E33").Value = "  +2.74%  "
$ws.Range("D34").Value = "593.62"
$ws.Range("E34").Value = "  +4.87%  "
$ws.Range("E35").Value = "  +5.00%  "
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("E37").Value = "  +5.39%  "
$ws.Range("D38").Value = "3.57"
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("E39").Value = "  +5.92%  "
$ws.Range("D40").Value = "36.35"
$ws.Range("E40").Value = "  +3.76%  "
$ws.Range("D41").Value = "0.0₃0768"
$ws.Range("E41").Value = "  +4.92%  "
$ws.Range("D42").Value = "3.195.58"
$ws.Range("E42").Value = "  +4.54%  "
$ws.Range("D43").Value = "2.93"
$ws.Range("E43").Value = "  +5.84%  "
$ws.Range("E44").Value = "  +4.49%  "
$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").Value = "2.55"
$ws.Range("E45").Value = "  +5.55%  "
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").Value = "2.79"
$ws.Range("E46").Value = "  +23.53%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "3.23"
$ws.Range("E47").Value = "  +0.79%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "0.135"
$ws.Range("E48").Value = "  +1.69%  "
$ws.Range("B49").Value = "FirstDigitalUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").Value = "8.71"
$ws.Range("E50").Value = "  +7.31%  "
$ws.Range("D51").Value = "140.00"
$ws.Range("E51").Value = "  -1.42%  "
